# Split the run " im Alter von ALTER Jahren, hier in STERBEORT."
# into three runs:
#   " im Alter von "
#   "LEBENS"
#   "ALTER Jahren, hier in STERBEORT."
# by typing "LEBENS" in front of "ALTER" and nudging the font of the
# inserted word so Word (COM) materialises it as its own run instead of
# silently re-merging it into the surrounding (identically formatted) text.

$d = $word.ActiveDocument

$target = $d.Content
$found = $target.Find.Execute("ALTER Jahren, hier in STERBEORT.")

if (-not $found) {
    throw "Could not locate 'ALTER Jahren, hier in STERBEORT.' in the document."
}

$insertStart = $target.Start

# Type the new word right before "ALTER Jahren, hier in STERBEORT."
$target.InsertBefore("LEBENS")

# Range covering the freshly typed "LEBENS" text.
$newWord = $d.Range($insertStart, $insertStart + 6)

# Re-stamp the (already correct) Arial font on just the new word so the
# run-splitting logic breaks it into its own <w:r>, keeping rFonts intact
# across ascii/hAnsi/complex-script.
$newWord.Font.Name = "Arial"
$newWord.Font.NameOther = "Arial"
$newWord.Font.NameBi = "Arial"

Write-Output $d.Content.Text
